$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 711.52856
$ws.Range("I15").Value2 = 711.52856
$ws.Range("K15").Value2 = 2134.58568
$ws.Range("M15").Value2 = -1965.58568

$ws.Range("H96").Value2 = 1385.2727
$ws.Range("I96").Value2 = 1619.8572
$ws.Range("K96").Value2 = 4859.571599999999
$ws.Range("M96").Value2 = -3486.571599999999

$ws.Range("H107").Value2 = 592.8
$ws.Range("I107").Value2 = 652.5238000000001
$ws.Range("J107").Value2 = 279.25
$ws.Range("K107").Value2 = 652.5238000000001
$ws.Range("L107").Value2 = 279.25
$ws.Range("M107").Value2 = 1267.4762
$ws.Range("N107").Value2 = -4119.25

$ws.Range("H113").Value2 = 3482.5833
$ws.Range("I113").Value2 = 1758.4
$ws.Range("J113").Value2 = 4714.143
$ws.Range("K113").Value2 = 1758.4
$ws.Range("L113").Value2 = 4714.143
$ws.Range("M113").Value2 = 1495.6
$ws.Range("N113").Value2 = -11222.143

$ws.Range("H138").Value2 = 380919.28
$ws.Range("I138").Value2 = 72642.86
$ws.Range("J138").Value2 = 668643.9399999999
$ws.Range("K138").Value2 = 217928.58
$ws.Range("L138").Value2 = 2005931.82
$ws.Range("M138").Value2 = -212788.58
$ws.Range("N138").Value2 = -2016211.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1146.8462
$ws.Range("I2").Value2 = 1148.05
$ws.Range("J2").Value2 = 1142.8334
$ws.Range("K2").Value2 = 1148.05
$ws.Range("L2").Value2 = 1142.8334
$ws.Range("M2").Value2 = -1035.05
$ws.Range("N2").Value2 = -1368.8334

$ws.Range("H45").Value2 = 2356.182
$ws.Range("I45").Value2 = 2067.3076
$ws.Range("K45").Value2 = 2067.3076
$ws.Range("M45").Value2 = -1690.3076

$ws.Range("H102").Value2 = 4407.643
$ws.Range("I102").Value2 = 2260.4443
$ws.Range("J102").Value2 = 8272.6
$ws.Range("K102").Value2 = 2260.4443
$ws.Range("L102").Value2 = 8272.6
$ws.Range("M102").Value2 = -638.4443000000001
$ws.Range("N102").Value2 = -11516.6

$ws.Range("H116").Value2 = 1146.8462
$ws.Range("I116").Value2 = 1148.05
$ws.Range("J116").Value2 = 1142.8334
$ws.Range("K116").Value2 = 1148.05
$ws.Range("L116").Value2 = 1142.8334
$ws.Range("M116").Value2 = 1145.95
$ws.Range("N116").Value2 = -5730.8334

$ws.Range("H122").Value2 = 4154
$ws.Range("I122").Value2 = 3668.7144
$ws.Range("K122").Value2 = 11006.1432
$ws.Range("M122").Value2 = -8556.143199999999

$ws.Range("H132").Value2 = 1902.5834
$ws.Range("I132").Value2 = 1733.683
$ws.Range("J132").Value2 = 2891.8572
$ws.Range("K132").Value2 = 5201.049
$ws.Range("L132").Value2 = 8675.571599999999
$ws.Range("M132").Value2 = -2671.049
$ws.Range("N132").Value2 = -13735.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1146.8462
$ws.Range("I3").Value2 = 1148.05
$ws.Range("J3").Value2 = 1142.8334
$ws.Range("K3").Value2 = 1148.05
$ws.Range("L3").Value2 = 1142.8334
$ws.Range("M3").Value2 = -1034.05
$ws.Range("N3").Value2 = -1370.8334

$ws.Range("H26").Value2 = 37249.5
$ws.Range("I26").Value2 = 37249.5
$ws.Range("K26").Value2 = 37249.5
$ws.Range("M26").Value2 = -36957.5

$ws.Range("H86").Value2 = 8511.823
$ws.Range("J86").Value2 = 2462.8
$ws.Range("L86").Value2 = 2462.8
$ws.Range("N86").Value2 = -4708.8

$ws.Range("H89").Value2 = 8511.823
$ws.Range("J89").Value2 = 2462.8
$ws.Range("L89").Value2 = 12314
$ws.Range("N89").Value2 = -23546

$ws.Range("H107").Value2 = 1737.875
$ws.Range("I107").Value2 = 1484
$ws.Range("K107").Value2 = 1484
$ws.Range("M107").Value2 = 436

$ws.Range("H134").Value2 = 6404.982
$ws.Range("I134").Value2 = 1407.3158
$ws.Range("J134").Value2 = 9042.638999999999
$ws.Range("K134").Value2 = 4221.9474
$ws.Range("L134").Value2 = 27127.917
$ws.Range("M134").Value2 = -1686.9474
$ws.Range("N134").Value2 = -32197.917

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2068.3713
$ws.Range("I31").Value2 = 1436.1818
$ws.Range("J31").Value2 = 12499.5
$ws.Range("K31").Value2 = 1436.1818
$ws.Range("L31").Value2 = 12499.5
$ws.Range("M31").Value2 = -1141.1818
$ws.Range("N31").Value2 = -13089.5

$ws.Range("H34").Value2 = 2068.3713
$ws.Range("I34").Value2 = 1436.1818
$ws.Range("J34").Value2 = 12499.5
$ws.Range("K34").Value2 = 1436.1818
$ws.Range("L34").Value2 = 12499.5
$ws.Range("M34").Value2 = -1234.1818
$ws.Range("N34").Value2 = -12903.5

$ws.Range("H58").Value2 = 2232.5
$ws.Range("I58").Value2 = 2400.5715
$ws.Range("K58").Value2 = 2400.5715
$ws.Range("M58").Value2 = -2197.5715

$ws.Range("H86").Value2 = 6700.5386
$ws.Range("I86").Value2 = 5142.4287
$ws.Range("K86").Value2 = 5142.4287
$ws.Range("M86").Value2 = -4019.4287

$ws.Range("H89").Value2 = 6700.5386
$ws.Range("I89").Value2 = 5142.4287
$ws.Range("K89").Value2 = 25712.1435
$ws.Range("M89").Value2 = -20096.1435

$ws.Range("H105").Value2 = 2405.5
$ws.Range("I105").Value2 = 2405.5
$ws.Range("K105").Value2 = 2405.5
$ws.Range("M105").Value2 = -658.5

$ws.Range("H132").Value2 = 4449911.5
$ws.Range("I132").Value2 = 5720529.5
$ws.Range("J132").Value2 = 2750
$ws.Range("K132").Value2 = 17161588.5
$ws.Range("L132").Value2 = 8250
$ws.Range("M132").Value2 = -17159058.5
$ws.Range("N132").Value2 = -13310

$ws.Range("H136").Value2 = 2232.5
$ws.Range("I136").Value2 = 2400.5715
$ws.Range("K136").Value2 = 7201.7145
$ws.Range("M136").Value2 = -4651.7145

$ws.Range("H141").Value2 = 195403.88
$ws.Range("J141").Value2 = 195403.88
$ws.Range("L141").Value2 = 195403.88
$ws.Range("N141").Value2 = -205763.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value2 = 1520
$ws.Range("I19").Value2 = 270
$ws.Range("J19").Value2 = 4020
$ws.Range("K19").Value2 = 810
$ws.Range("L19").Value2 = 12060
$ws.Range("M19").Value2 = -636
$ws.Range("N19").Value2 = -12408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 31759.766
$ws.Range("I102").Value2 = 2076.76
$ws.Range("K102").Value2 = 2076.76
$ws.Range("M102").Value2 = -454.7600000000002

$ws.Range("H122").Value2 = 4477.231
$ws.Range("I122").Value2 = 4175.625
$ws.Range("K122").Value2 = 12526.875
$ws.Range("M122").Value2 = -10076.875

$ws.Range("H132").Value2 = 11499059
$ws.Range("I132").Value2 = 17548512
$ws.Range("K132").Value2 = 52645536
$ws.Range("M132").Value2 = -52643006

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value2 = 38995
$ws.Range("J25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("N25").ClearContents()

$ws.Range("H40").Value2 = 4356
$ws.Range("I40").Value2 = 4346
$ws.Range("J40").Value2 = 4381.7144
$ws.Range("K40").Value2 = 4346
$ws.Range("L40").Value2 = 4381.7144
$ws.Range("M40").Value2 = -4210
$ws.Range("N40").Value2 = -4653.7144

$ws.Range("H61").Value2 = 13968.75
$ws.Range("I61").Value2 = 17791.666
$ws.Range("K61").Value2 = 17791.666
$ws.Range("M61").Value2 = -17589.666

$ws.Range("H93").Value2 = 4687.778
$ws.Range("I93").Value2 = 4523.75
$ws.Range("K93").Value2 = 4523.75
$ws.Range("M93").Value2 = -3275.75

$ws.Range("H113").Value2 = 13968.75
$ws.Range("I113").Value2 = 17791.666
$ws.Range("K113").Value2 = 17791.666
$ws.Range("M113").Value2 = -15621.666

$ws.Range("H132").Value2 = 2885.2952
$ws.Range("I132").Value2 = 2859.8164
$ws.Range("J132").Value2 = 2989.3333
$ws.Range("K132").Value2 = 8579.449200000001
$ws.Range("L132").Value2 = 8967.999899999999
$ws.Range("M132").Value2 = -6049.449200000001
$ws.Range("N132").Value2 = -14027.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value2 = 24997
$ws.Range("I40").Value2 = 19995
$ws.Range("J40").Value2 = 29999
$ws.Range("K40").Value2 = 19995
$ws.Range("L40").Value2 = 29999
$ws.Range("M40").Value2 = -19846
$ws.Range("N40").Value2 = -30297

$ws.Range("H122").Value2 = 9998.75
$ws.Range("I122").Value2 = 9998.333000000001
$ws.Range("K122").Value2 = 29994.999
$ws.Range("M122").Value2 = -27544.999

$ws.Range("H132").Value2 = 3386.963
$ws.Range("I132").Value2 = 3732.3
$ws.Range("K132").Value2 = 11196.9
$ws.Range("M132").Value2 = -8666.900000000001

$ws.Range("H136").Value2 = 9961.305
$ws.Range("I136").Value2 = 18184.334
$ws.Range("K136").Value2 = 54553.00199999999
$ws.Range("M136").Value2 = -52003.00199999999
